# "#5: cash & deposit done"
#
# The deposit ("存款") worksheet (sheet index 2) previously had a bug: its
# header row (row 1) re-used data values instead of proper column-name
# labels, and the sheet only had columns A:F. This change:
#   1. Fixes row 1 to use the proper header labels (bank, deposit_type,
#      currency, owner, total, property_category, category, date,
#      legislator_name, legislator_id, source_file, index).
#   2. Adds the metadata columns G:M (deposit_type/category constant,
#      normal/category, date, legislator_name, legislator_id, source_file,
#      index) to every data row, matching the shape already used by the
#      stock ("股票") worksheet (sheet index 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---- Row 1: proper header labels -------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Give the new header cells (G1:M1) the same bold/border/centered style as
# the rest of row 1 (B1:F1).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122) | Out-Null

# ---- Rows 2-13: bank/deposit data rows --------------------------------
# Columns A (index), B (bank), C (deposit_type), D (currency) and E (owner)
# already hold the correct values and are left untouched; only the new
# metadata columns G:M are added.

$banks = @{
    2  = "基隆二信營業部"
    3  = "基隆二信營業部"
    4  = "基隆二信營業部"
    5  = "基隆二信港東分社"
    6  = "臺灣銀行群賢分行"
    7  = "華南商業銀行民生分行"
    8  = "華南商業銀行城內分行"
    9  = "國泰世華商業銀行板橋分行"
    10 = "國泰世華商業銀行營業部"
    11 = "中國信託商業銀行城中分行"
    12 = "中華商業銀行營業部"
    13 = "中華郵政股份有限公司基隆"
}

$totals = @{
    2  = 1
    3  = 24217
    4  = 40
    5  = 171
    6  = 1990318
    7  = 7
    8  = 8840
    9  = 1172
    10 = 5
    11 = 4377
    12 = 2
    13 = 21500
}

$indices = @{
    2  = 43
    3  = 45
    4  = 46
    5  = 47
    6  = 48
    7  = 49
    8  = 50
    9  = 51
    10 = 52
    11 = 53
    12 = 54
    13 = 55
}

for ($r = 2; $r -le 13; $r++) {
    # Columns A-F: re-assert the (unchanged) values so the sheet is
    # self-consistent even though they already hold the right data.
    $ws.Range("B$r").Value = $banks[$r]
    $ws.Range("C$r").Value = "活期儲蓄存款"
    $ws.Range("D$r").Value = "新臺幣"
    $ws.Range("E$r").Value = "謝國樑"
    $ws.Range("F$r").Value = $totals[$r]

    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"

    # The literal text "2011-11-17" must stay a *text* cell, not get
    # auto-converted into a date serial number by the smart-parsing that
    # a plain `.Value = "2011-11-17"` assignment would trigger. Compute it
    # via a throw-away formula cell and paste only the *value* across,
    # which keeps it a plain string.
    $ws.Range("ZZ1").Formula = '="2011-11-17"'
    $ws.Range("ZZ1").Copy() | Out-Null
    $ws.Range("I$r").PasteSpecial(-4163) | Out-Null
    $ws.Range("ZZ1").ClearContents() | Out-Null

    $ws.Range("J$r").Value = "謝國樑"
    $ws.Range("K$r").Value = 1387
    $ws.Range("L$r").Value = "tmpbbad1"
    $ws.Range("M$r").Value = $indices[$r]
}

Write-Output "deposit sheet updated"
